$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TFX")

# Row 4 - Inventory
$ws.Range("B4").Value = 512000000.0
$ws.Range("C4").Value = 513000000.0
$ws.Range("D4").Value = 526000000.0
$ws.Range("E4").Value = 515000000.0
$ws.Range("F4").Value = 489000000.0

# Row 15 - Accounts Payable
$ws.Range("C15").Value = 103000000.0
$ws.Range("D15").Value = 97000000.0
$ws.Range("E15").Value = 107000000.0
$ws.Range("F15").Value = 104000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 474000000.0
$ws.Range("C24").Value = 477000000.0
$ws.Range("D24").Value = 481000000.0
$ws.Range("E24").Value = 483000000.0
$ws.Range("F24").Value = 484000000.0

# Row 37 - Net Debt
$ws.Range("G37").Value = 1607860000.0

# Row 38 - Total Debt
$ws.Range("G38").Value = 1908943000.0
